$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.098.04"
$ws.Cells.Item(2, 5).Value = "  -0.06%  "

$ws.Cells.Item(3, 4).Value = "1.789.54"
$ws.Cells.Item(3, 5).Value = "  +0.21%  "

$ws.Cells.Item(4, 5).Value = "  +0.10%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "227.57"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.43%  "

$ws.Cells.Item(6, 5).Value = "  -0.56%  "

$ws.Cells.Item(7, 5).Value = "  +0.10%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "32.24"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.47%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.293"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.84%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0692"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.88%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0940"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.31%  "

$ws.Cells.Item(12, 4).Value = "2.047.63"
$ws.Cells.Item(12, 5).Value = "  +0.22%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "11.56"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +4.50%  "

$ws.Cells.Item(14, 4).Value = "1.782.20"
$ws.Cells.Item(14, 5).Value = "  -0.07%  "

$ws.Cells.Item(15, 2).Value = "WrappedBTC"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(15, 4).Value = "34.107.07"
$ws.Cells.Item(15, 5).Value = "  +0.18%  "

$ws.Cells.Item(16, 2).Value = "Polygon"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.622"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.06%  "

$ws.Cells.Item(17, 5).Value = "  +0.48%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "67.83"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.02%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "245.11"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.12%  "

$ws.Cells.Item(20, 5).Value = "  -1.00%  "

$ws.Cells.Item(21, 5).Value = "  +0.07%  "

$ws.Cells.Item(22, 5).Value = "  +0.21%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.10"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.25%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.05"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.59%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "162.10"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.74%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.16"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.32%  "

$ws.Cells.Item(27, 5).Value = "  -0.17%  "

$ws.Cells.Item(28, 5).Value = "  +0.83%  "

$ws.Cells.Item(29, 5).Value = "  +0.25%  "

$ws.Cells.Item(30, 5).Value = "  +2.12%  "

$ws.Cells.Item(31, 5).Value = "  +1.50%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.66"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.26%  "

$ws.Cells.Item(33, 5).Value = "  +2.92%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.83"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.59%  "

$ws.Cells.Item(35, 4).Value = "1.421.02"
$ws.Cells.Item(35, 5).Value = "  +1.80%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.642"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.31%  "

$ws.Cells.Item(37, 5).Value = "  +2.31%  "

$ws.Cells.Item(38, 5).Value = "  +7.10%  "

$ws.Cells.Item(39, 5).Value = "  -0.71%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "80.54"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.10%  "

$ws.Cells.Item(41, 5).Value = "  +0.63%  "

$ws.Cells.Item(42, 5).Value = "  +0.60%  "

$ws.Cells.Item(43, 5).Value = "  -0.12%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "13.33"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +7.57%  "

$ws.Cells.Item(45, 5).Value = "  +3.52%  "

$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(47, 4).Value = "0.0₆0137"
$ws.Cells.Item(47, 5).Value = "  -2.14%  "

$ws.Cells.Item(48, 2).Value = "WEMIXToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.07"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.13%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "107.52"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.63%  "

$ws.Cells.Item(50, 4).Value = "1.948.80"
$ws.Cells.Item(50, 5).Value = "  +0.24%  "

$ws.Cells.Item(51, 5).Value = "  +0.10%  "

